# Apply changes described in the commit:
# "Tasks + Bucket's Attachments + Update RD&T"
#
# - Update Bucket/Attachments sub-task values (rows 26-29 feed D25)
# - Update RD&T-related sub-task values (rows 42-49 feed D41 and D46)
# - Formulas in D25, D41, D46 will recalculate automatically
# - Update the sheet view's scroll position / selection to match final state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Bucket section - Attachments task (row 28) bumped from 50 to 100
$ws.Range("D28").Value = 100

# RD&T section - tasks at rows 42-45 (feeding D41) updated
$ws.Range("D42").Value = 100
$ws.Range("D43").Value = 100
$ws.Range("D44").Value = 100
$ws.Range("D45").Value = 50

# RD&T section - tasks at rows 47-49 (feeding D46) updated
$ws.Range("D48").Value = 60

# Force recalculation so dependent SUM formulas (D25, D41, D46) refresh
$excel.CalculateFullRebuild()
$wb.Application.Calculate()

# Update the view state to match what was saved (scrolled down / new selection)
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D49").Select()
